$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy header style (bold, border, centered) from R1 into S1:T1
$ws.Range("R1").Copy() | Out-Null
$ws.Range("S1:T1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# Header row values
$ws.Range("S1").Value = 'test7'
$ws.Range("T1").Value = 'emotion_test7'

# Data rows 2-22: column S (response text) and column T (emotion label)
$ws.Range("S2").Value = 'Entiendo cómo te sientes. El estrés antes de un examen es algo común. Una buena idea es repasar lo más importante y tomar descansos cortos para no saturarte. También puedes practicar la respiración profunda antes de empezar, eso te ayudará a sentirte más tranquilo. Recuerda que has trabajado duro y estás listo para enfrentar el examen. ¡Tú puedes! 💪'
$ws.Range("T2").Value = 'sadness'
$ws.Range("S3").Value = '¡Qué bien! Felicitaciones por esa buena nota en la expo. Todo el esfuerzo y las trasnochadas valieron la pena. Esos momentos de sacrificio a menudo traen grandes recompensas. ¿Tienes algún plan para celebrar tu éxito? ¡Disfrútalo! 🎉✨'
$ws.Range("T3").Value = 'joy'
$ws.Range("S4").Value = 'Lamento que no hayas podido descansar bien. Es normal sentirse agotado después de una noche de estudio. Te recomendaría que tomes un pequeño descanso hoy, si puedes, y que te hidrates. También podrías intentar hacer una pausa para respirar conscientemente durante un minuto; solo observa tu respiración, eso puede ayudarte a relajarte un poco. Recuerda que es importante cuidar de ti mismo, así que trata de dormir bien esta noche. ¡Ánimo!'
$ws.Range("T4").Value = 'sadness'
$ws.Range("S5").Value = '¡Qué bien que te sientas así! La buena energía es contagiosa y puede hacer que tu día sea mucho mejor. Aprovecha esa vibra positiva para hacer algo que te guste, ya sea estudiar, salir a caminar o simplemente disfrutar de un buen rato con amigos. Recuerda que esos momentos son los que hacen que la vida sea especial. ¡Sigue disfrutando de esa buena energía! ✨😊'
$ws.Range("T5").Value = 'joy'
$ws.Range("S6").Value = 'Es normal sentir miedo ante la posibilidad de no pasar una materia, pero recuerda que esos pensamientos pueden ser solo eso, pensamientos. En lugar de pensar en lo peor, enfócate en lo que puedes hacer ahora para mejorar. Establece un plan de estudio y busca ayuda si la necesitas. Recuérdate a ti mismo que estás haciendo lo mejor que puedes. ¡Tú puedes superar este desafío!'
$ws.Range("T6").Value = 'fear'
$ws.Range("S7").Value = 'Es normal sentir que el tiempo no alcanza, especialmente cuando tienes muchas cosas en la cabeza. Te recomendaría hacer una lista de tareas y priorizar lo más importante. A veces, dividir las tareas en partes más pequeñas puede hacer que se sientan más manejables. Recuerda que está bien pedir ayuda si la necesitas. ¡Tómalo un paso a la vez!'
$ws.Range("T7").Value = 'sadness'
$ws.Range("S8").Value = '¡Eso es genial! Hacer ejercicio es una excelente manera de liberar el estrés y mejorar tu estado de ánimo. Además, te ayuda a sentirte más enérgico y enfocado. Sigue así, ¡es un gran hábito! ¿Tienes algún tipo de ejercicio favorito que disfrutes hacer? 💪😊'
$ws.Range("T8").Value = 'joy'
$ws.Range("S9").Value = 'Es completamente normal sentir ansiedad al hablar en público o participar en clase. Muchos sienten lo mismo. Una técnica que puede ayudarte es practicar lo que quieres decir en voz alta antes de hacerlo; eso puede hacerte sentir más preparado. También, intenta recordar que todos están ahí para aprender, así que no eres el único que puede sentirse nervioso. ¡Con el tiempo y la práctica, te sentirás más cómodo! ¡Tú puedes!'
$ws.Range("T9").Value = 'fear'
$ws.Range("S10").Value = 'Entiendo cómo te sientes. A veces, la desmotivación puede ser abrumadora. Una buena idea es establecer pequeñas metas diarias; eso puede ayudarte a sentir un sentido de logro. También, intenta cambiar de ambiente o hacer algo diferente para refrescar tu mente. Recuerda que está bien tener días bajos, pero lo importante es seguir avanzando, aunque sea un pasito a la vez. ¡Tú puedes!'
$ws.Range("T10").Value = 'sadness'
$ws.Range("S11").Value = '¡Eso es increíble! 🎉 Felicitaciones por tu logro en el parcial de inteligencia artificial. Es genial ver que tu esfuerzo está dando frutos. Celebra ese éxito, te lo mereces. ¿Tienes algún plan para festejarlo? ¡Sigue así!'
$ws.Range("T11").Value = 'joy'
$ws.Range("S12").Value = 'Lamento que te sientas así. Es importante reconocer que el cansancio puede afectar mucho nuestra motivación. A veces, tomarse un pequeño descanso puede hacer maravillas. Permítete un momento para recargar energías, ya sea descansando, haciendo algo que disfrutes o simplemente desconectando un rato. Recuerda que está bien pedir ayuda si sientes que lo necesitas. No estás sola en esto. ¡Tú puedes!'
$ws.Range("T12").Value = 'sadness'
$ws.Range("S13").Value = 'Lamento que te sientas así. Es normal tener momentos de duda sobre si estás en el lugar correcto. A veces, puede ayudar hablar con alguien sobre tus sentimientos, ya sea un amigo, un profesor o un consejero. Recuerda que cada camino es único y que está bien tomarse el tiempo para descubrir lo que realmente te apasiona. No estás solo en esto, y siempre hay oportunidades para explorar diferentes opciones. ¡Tú vales mucho!'
$ws.Range("T13").Value = 'sadness'
$ws.Range("S14").Value = '¡Eso es genial! 🎉 Organizar tu tiempo y cumplir con tus entregas es un gran logro. Te mereces sentirte orgulloso de ti mismo. Esa sensación de haber completado todo a tiempo es increíble y te motivará para seguir así. ¿Tienes algún plan para celebrar tu éxito? ¡Sigue así!'
$ws.Range("T14").Value = 'joy'
$ws.Range("S15").Value = 'Esto que me cuentas es muy importante y me importa mucho tu bienestar. Necesitas ayuda inmediata. Por favor, contacta a la Universidad del Valle - Tuluá a través del correo serviciopsicologico.tulua@correounivalle.edu.co. Tu vida tiene valor. Por favor, contacta estos recursos AHORA. No estás solo/a.'
$ws.Range("T15").Value = 'sadness'
$ws.Range("S16").Value = 'Entiendo que te sientas así, a veces el semestre puede ser abrumador. ¿Qué tal si te tomas un momento para respirar y pensar en pequeños pasos que puedas dar? Tal vez podrías hacer una lista de las cosas que te están estresando y ver si hay alguna manera de manejarlas una a una. Recuerda que está bien pedir ayuda si la necesitas. ¡Tú puedes con esto!'
$ws.Range("T16").Value = 'anger'
$ws.Range("S17").Value = 'Las discusiones en grupos pueden ser realmente difíciles. Lo primero que podrías hacer es tomarte un momento para calmarte y reflexionar sobre lo que sucedió. Luego, intenta hablar con tus compañeros de manera abierta y honesta. A veces, expresar cómo te sientes y escuchar a los demás puede ayudar a encontrar un terreno común. Recuerda que todos están ahí para aprender y trabajar juntos. ¡Confía en que puedes resolverlo!'
$ws.Range("T17").Value = 'others'
$ws.Range("S18").Value = 'Es completamente normal extrañar a tu familia y tu hogar, especialmente cuando estás en un lugar nuevo. La soledad puede ser difícil de manejar, pero hay formas de sentirte más conectado. Intenta programar videollamadas con tu familia o amigos, o incluso busca grupos de apoyo en la universidad donde puedas conocer a otros que se sientan igual. También podrías explorar actividades o clubes que te interesen; eso puede ayudarte a hacer nuevos amigos y sentirte más en casa. ¡No estás solo en esto!'
$ws.Range("T18").Value = 'sadness'
$ws.Range("S19").Value = 'Es normal sentir esa tentación, especialmente cuando la presión es alta. Pero recuerda que copiarte no es la solución y puede traerte más problemas a largo plazo. En lugar de eso, intenta repasar lo que has estudiado y confía en tus capacidades. Si sientes que necesitas más apoyo, considera hablar con un profesor o un compañero para aclarar dudas. ¡Tú puedes hacerlo!'
$ws.Range("T19").Value = 'fear'
$ws.Range("S20").Value = '¡Eso es genial! Terminar un proyecto a tiempo es un gran logro y es totalmente válido sentirse orgulloso de ti mismo. Celebra ese éxito, aunque sea con algo pequeño. Recuerda que cada paso cuenta y que estás avanzando en tu camino. ¿Qué tal si escribes ese logro en un lugar visible para recordarlo y motivarte en el futuro? ¡Sigue así! 🎉'
$ws.Range("T20").Value = 'joy'
$ws.Range("S21").Value = '¡Eso es increíble! Hacer nuevos amigos es una de las mejores partes de la experiencia universitaria. Tener a alguien con quien compartir tus días puede hacer que todo sea más divertido y menos estresante. ¿Ya han hecho algo juntos o tienen planes para salir? ¡Aprovecha esa nueva amistad! 🎉'
$ws.Range("T21").Value = 'joy'
$ws.Range("S22").Value = 'Parece que no has escrito nada. Si necesitas hablar de algo o tienes alguna pregunta, aquí estoy para ayudarte. ¡Cuéntame!'
$ws.Range("T22").Value = 'others'
